$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that disappear from the table (old row 14: bitcoin/com.hamxa.shaynachim,
# and old row 16 (after first delete becomes row 15): blank A / com.sugar.powerfulquotes).
# Deleting entire rows shifts everything below up automatically.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(15).Delete()

# Rows 14-15 now already hold the correct final content (blockchain technology / bitcoin rows).
# Rows 16-19 need their content + height rewritten to match the new table below them.
$ws.Rows.Item(16).RowHeight = 24
$ws.Cells.Item(16, 1).Value = "affiliate marketing"
$ws.Cells.Item(16, 2).Value = "affiliate.marketing.guide"

$ws.Rows.Item(17).RowHeight = 46.5
$ws.Cells.Item(17, 1).Value = "Powerful Positive Motivation Quotes"
$ws.Cells.Item(17, 2).Value = "com.sugar.powerfulquotes"

$ws.Rows.Item(18).RowHeight = 12.8
$ws.Cells.Item(18, 1).Value = "affiliate marketing"
$ws.Cells.Item(18, 2).Value = "affiliate.marketing.guide"

$ws.Rows.Item(19).RowHeight = 12.8
$ws.Cells.Item(19, 1).Value = "passive income"
$ws.Cells.Item(19, 2).Value = "passive.income.nadi.myfirstdrawermenuproject2"

# Add a new trailing row 20, cloning formatting from row 19 so the style index matches
# the rest of the table, then set its values.
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 12.8
$ws.Cells.Item(20, 1).Value = "bitcoin"
$ws.Cells.Item(20, 2).Value = "com.hamxa.shaynachim"

# Match the updated selection recorded in the saved file (engine does not model/serialize
# the sheetView's topLeftCell scroll position, only the active cell/selection).
$ws.Range("A15").Select()
